$wb = $excel.ActiveWorkbook

# Update the saved selection on the "raw" sheet (view-state only, no data change)
$wsRaw = $wb.Worksheets.Item("raw")
$wsRaw.Activate()
$wsRaw.Range("J34").Select()

$ws = $wb.Worksheets.Item("2000-2011")

# Make this the active sheet (matches after-state activeTab/tabSelected)
$ws.Activate()

# Insert a new row at position 18, shifting existing rows (18-23) down to (19-24)
$ws.Rows.Item(18).Insert()

# Copy formatting from row 17 (the row above) onto the newly inserted row 18
$ws.Range("A17:V17").Copy()
$ws.Range("A18:V18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The paste-formats step also stamps a few columns that have no formatting
# of their own in row 17 (and shouldn't exist as explicit cells) - drop them
$ws.Range("Q18").Clear()
$ws.Range("T18").Clear()
$ws.Range("V18").Clear()

# Fill in the new benchmark entry: "Testbench Ideacenter" w/ Pentium E6600 + nVidia 9600 GSO
$ws.Cells.Item(18, 1).Value = "Testbench Ideacenter"
$ws.Cells.Item(18, 2).Value = "Pentium E6600"
$ws.Cells.Item(18, 3).Value = 3.06
$ws.Cells.Item(18, 4).Value = "nVidia 9600 GSO"
$ws.Cells.Item(18, 12).Value = 8052
$ws.Cells.Item(18, 13).Value = 3753
$ws.Cells.Item(18, 14).Value = 3625
$ws.Cells.Item(18, 15).Value = 1873

# Update the selection on this sheet to match the saved view state
$ws.Range("G29").Select()
